# Fill in the missing low-intake (0 .. boundary) relative-risk (RR) values for
# several food groups in the RR interpolation table.
#
# Background: for every food group row, column B holds the RR at intake = 0
# (the header row stores the intake amount, with column B = intake 0, column
# C = intake 1, etc.). RR at intake = 0 must always be 1 (no exposure -> no
# relative risk change), and the cells between intake = 0 and the first
# already-populated intake level were missing for a handful of rows. This
# script rebuilds those missing cells by linearly interpolating between
# RR = 1 at intake = 0 and the RR already stored at the first populated
# intake level for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> the column letter of the first already-populated cell
# (the right-hand boundary of the interpolation) for that row.
$rowBoundaryCol = @{
    3  = "CX"
    4  = "AA"
    5  = "AZ"
    7  = "CX"
    8  = "AA"
    11 = "AZ"
    13 = "AZ"
    15 = "AZ"
}

function ColLettersToNumber([string]$colLetters) {
    $n = 0
    foreach ($ch in $colLetters.ToCharArray()) {
        $n = $n * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $n
}

function NumberToColLetters([int]$n) {
    $s = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $s = [char](65 + $rem) + $s
        $n = [int](($n - 1) / 26)
    }
    return $s
}

foreach ($rowNum in $rowBoundaryCol.Keys) {
    $boundaryColLetters = $rowBoundaryCol[$rowNum]
    $boundaryColNum = ColLettersToNumber $boundaryColLetters

    # "Intake" x-position of the boundary column: column B (=2) is intake 0,
    # so intake = columnNumber - 2.
    $xBoundary = $boundaryColNum - 2

    # RR value already present at the boundary column (right edge of the gap).
    $yBoundary = $ws.Cells.Item($rowNum, $boundaryColNum).Value2
    $yStart = 1.0

    # Number of cells to fill: columns B .. (boundary - 1)
    $fillCount = $xBoundary
    if ($fillCount -le 0) { continue }

    $lastFillColNum = $boundaryColNum - 1
    $lastFillColLetters = NumberToColLetters $lastFillColNum

    $values = New-Object 'object[,]' 1,$fillCount
    for ($i = 0; $i -lt $fillCount; $i++) {
        $x = $i
        $values[0, $i] = $yStart + ($yBoundary - $yStart) * $x / $xBoundary
    }

    $rangeAddress = "B" + $rowNum + ":" + $lastFillColLetters + $rowNum
    $ws.Range($rangeAddress).Value = $values
}
